$d = $word.ActiveDocument

$replacements = @(
    @("2025-05-07 Wednesday", "2025-05-08 Thursday"),
    @("191÷3=63, 2", "674÷7=96, 2"),
    @("901÷8=112, 5", "410÷4=102, 2"),
    @("319÷8=39, 7", "745÷4=186, 1"),
    @("432÷8=54, 0", "626÷8=78, 2"),
    @("841÷5=168, 1", "453÷6=75, 3"),
    @("771÷2=385, 1", "337÷8=42, 1"),
    @("771÷3=257, 0", "217÷4=54, 1"),
    @("800÷9=88, 8", "434÷2=217, 0"),
    @("932÷5=186, 2", "220÷8=27, 4"),
    @("639÷9=71, 0", "539÷3=179, 2"),
    @("640÷5=128, 0", "254÷6=42, 2"),
    @("614÷7=87, 5", "822÷6=137, 0"),
    @("702÷7=100, 2", "400÷3=133, 1"),
    @("394÷9=43, 7", "999÷3=333, 0"),
    @("112÷2=56, 0", "598÷5=119, 3"),
    @("758÷3=252, 2", "364÷3=121, 1"),
    @("876÷9=97, 3", "593÷7=84, 5"),
    @("449÷6=74, 5", "681÷2=340, 1"),
    @("340÷7=48, 4", "159÷3=53, 0"),
    @("553÷8=69, 1", "292÷4=73, 0"),
    @("510÷3=170, 0", "246÷2=123, 0"),
    @("777÷5=155, 2", "561÷3=187, 0"),
    @("805÷9=89, 4", "427÷8=53, 3"),
    @("717÷6=119, 3", "974÷4=243, 2"),
    @("828÷8=103, 4", "248÷5=49, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
